$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F header: "VLJ ID" -> "CSS ID"
$ws.Range("F1").Value = "CSS ID"
# Nudge the font so a distinct (but visually identical, black) style record
# is created for the header cell, matching the extra cellXf/font produced
# when Excel re-typed this header.
$ws.Range("F1").Font.ThemeColor = 1

# Column F row 2: "123" -> "BVAHUELS"
$ws.Range("F2").Value = "BVAHUELS"

# Row 3 column F ("456") and column G (VLJ / Huels, Stuart) are unchanged.

# Leave the active selection on G6, matching the saved workbook view.
$ws.Range("G6").Select() | Out-Null

Write-Host "edit complete"
